# Rename the "Layman " sheet (trailing space) to "Layman"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Layman ")
$ws.Name = "Layman"

# Make "Simple Layman" the active/selected sheet (was "PrePreschool")
$active = $wb.Worksheets.Item("Simple Layman")
$active.Activate()
